$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 1017
$ws.Range("I45").Value = 1017
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 3051
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -2859
$ws.Range("N45").ClearContents()

$ws.Range("H80").Value = 1920
$ws.Range("I80").Value = 2500
$ws.Range("J80").Value = 1533.3334
$ws.Range("K80").Value = 7500
$ws.Range("L80").Value = 4600.0002
$ws.Range("M80").Value = -6502
$ws.Range("N80").Value = -6596.0002

$ws.Range("H83").Value = 1920
$ws.Range("I83").Value = 2500
$ws.Range("J83").Value = 1533.3334
$ws.Range("K83").Value = 22500
$ws.Range("L83").Value = 13800.0006
$ws.Range("M83").Value = -17508
$ws.Range("N83").Value = -23784.0006

$ws.Range("H88").Value = 2666.3333
$ws.Range("I88").Value = 4999
$ws.Range("J88").Value = 1500
$ws.Range("K88").Value = 4999
$ws.Range("L88").Value = 1500
$ws.Range("M88").Value = -4593
$ws.Range("N88").Value = -2312

$ws.Range("H91").Value = 2666.3333
$ws.Range("I91").Value = 4999
$ws.Range("J91").Value = 1500
$ws.Range("K91").Value = 4999
$ws.Range("L91").Value = 1500
$ws.Range("M91").Value = -3595
$ws.Range("N91").Value = -4308

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 199.5
$ws.Range("I10").Value = 199.5
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 199.5
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -59.5
$ws.Range("N10").ClearContents()

$ws.Range("H86").Value = 6183.3335
$ws.Range("I86").Value = 6183.3335
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 6183.3335
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -5060.3335

$ws.Range("H89").Value = 6183.3335
$ws.Range("I89").Value = 6183.3335
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 30916.6675
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -25300.6675

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2965
$ws.Range("I31").Value = 1948.3334
$ws.Range("J31").Value = 3519.5454
$ws.Range("K31").Value = 1948.3334
$ws.Range("L31").Value = 3519.5454
$ws.Range("M31").Value = -1653.3334
$ws.Range("N31").Value = -4109.5454

$ws.Range("H34").Value = 2965
$ws.Range("I34").Value = 1948.3334
$ws.Range("J34").Value = 3519.5454
$ws.Range("K34").Value = 1948.3334
$ws.Range("L34").Value = 3519.5454
$ws.Range("M34").Value = -1746.3334
$ws.Range("N34").Value = -3923.5454

$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -550

$ws.Range("H132").Value = 3133.8333
$ws.Range("I132").Value = 3252.2307
$ws.Range("J132").Value = 2826
$ws.Range("K132").Value = 9756.6921
$ws.Range("L132").Value = 8478
$ws.Range("M132").Value = -7226.6921
$ws.Range("N132").Value = -13538

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3500
$ws.Range("I3").Value = 3500
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 10500
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -10388

$ws.Range("H38").Value = 1327.1666
$ws.Range("I38").Value = 2543
$ws.Range("J38").Value = 111.333336
$ws.Range("K38").Value = 7629
$ws.Range("L38").Value = 334.000008
$ws.Range("M38").Value = -7282
$ws.Range("N38").Value = -1028.000008

$ws.Range("H131").Value = 2838.182
$ws.Range("I131").Value = 1244
$ws.Range("J131").Value = 4166.6665
$ws.Range("K131").Value = 3732
$ws.Range("L131").Value = 12499.9995
$ws.Range("M131").Value = 1308
$ws.Range("N131").Value = -22579.9995

$ws.Range("H134").Value = 3460
$ws.Range("I134").Value = 2166.6667
$ws.Range("J134").Value = 5400
$ws.Range("K134").Value = 6500.000100000001
$ws.Range("L134").Value = 16200
$ws.Range("M134").Value = -1430.000100000001
$ws.Range("N134").Value = -26340

$ws.Range("H141").Value = 3200
$ws.Range("I141").Value = 3200
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 9600
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -4420

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 94.2
$ws.Range("I2").Value = 82.333336
$ws.Range("J2").Value = 112
$ws.Range("K2").Value = 82.333336
$ws.Range("L2").Value = 112
$ws.Range("M2").Value = 30.666664
$ws.Range("N2").Value = -338

$ws.Range("H70").Value = 4750
$ws.Range("I70").Value = 4750
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 4750
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -4480

$ws.Range("H73").Value = 4750
$ws.Range("I73").Value = 4750
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 4750
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -3814

$ws.Range("H80").Value = 13507.143
$ws.Range("I80").Value = 2900
$ws.Range("J80").Value = 21462.5
$ws.Range("K80").Value = 2900
$ws.Range("L80").Value = 21462.5
$ws.Range("M80").Value = -1902
$ws.Range("N80").Value = -23458.5

$ws.Range("H83").Value = 13507.143
$ws.Range("I83").Value = 2900
$ws.Range("J83").Value = 21462.5
$ws.Range("K83").Value = 14500
$ws.Range("L83").Value = 107312.5
$ws.Range("M83").Value = -9508
$ws.Range("N83").Value = -117296.5

$ws.Range("H102").Value = 1102.4
$ws.Range("I102").Value = 1102.4
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1102.4
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 519.5999999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3999
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 3999
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 3999
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -4271

$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").ClearContents()

$ws.Range("H122").Value = 4664.2856
$ws.Range("I122").Value = 4691.6665
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 14074.9995
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -11624.9995
$ws.Range("N122").Value = -18400

$ws.Range("H132").Value = 1189.3334
$ws.Range("I132").Value = 1189.3334
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3568.0002
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1038.0002

$ws.Range("H141").Value = 99999.336
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 99999.336
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 99999.336
$ws.Range("N141").Value = -110359.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7652.75
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 7652.75
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 7652.75
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -8900.75

$ws.Range("H65").Value = 7652.75
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 7652.75
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 38263.75
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -44503.75

$ws.Range("H119").Value = 10000
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 10000
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 10000
$ws.Range("N119").Value = -19676

$ws.Range("H122").Value = 1750
$ws.Range("I122").Value = 1750
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5250
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2800

$ws.Range("H132").Value = 1338.7
$ws.Range("I132").Value = 1338.7
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4016.1
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1486.1
